$wb = $excel.ActiveWorkbook

# --- Update selection on sheet "Register" (was F14 -> now G14) ---
$wsRegister = $wb.Worksheets.Item("Register")
$wsRegister.Range("G14").Select()

# --- Add the two new worksheets at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsOrder = $wb.Worksheets.Add($null, $lastSheet)
$wsOrder.Name = "CommonForOrder"
$wsOrder.Range("A1").Value = "TaxForProducts"
$wsOrder.Range("A2").Value = 0.1
$wsOrder.Columns.Item(1).ColumnWidth = 14.7109375
$wsOrder.PageSetup.PaperSize = 9
$wsOrder.PageSetup.Orientation = 1
$wsOrder.Range("A3").Select()

$wsValidation = $wb.Worksheets.Add($null, $wsOrder)
$wsValidation.Name = "Step 2 validation"
$wsValidation.Range("B2:F4").Select()

# --- The newly added "CommonForOrder" sheet should be the active tab ---
$wsOrder.Select()
